$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume/Change (E) columns for the cryptos list refresh.
# The Price column holds values that look numeric (e.g. "605.14") but must
# remain plain text exactly as scraped (matching the original inline-string
# cells), so a leading apostrophe is used to force text entry and avoid
# Excel's automatic numeric conversion / precision or trailing-zero loss.

# Row 2 - Bitcoin
$ws.Range("D2").Value = "'69.862.41"
$ws.Range("E2").Value = "  +0.72%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "'3.530.16"
$ws.Range("E3").Value = "  +1.14%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.14%  "

# Row 5 - BNB
$ws.Range("D5").Value = "'605.14"
$ws.Range("E5").Value = "  -1.05%  "

# Row 6 - Solana
$ws.Range("D6").Value = "'197.09"
$ws.Range("E6").Value = "  +6.07%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +0.49%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.05%  "

# Row 9 - Dogecoin
$ws.Range("E9").Value = "  -5.44%  "

# Row 10 - Cardano
$ws.Range("E10").Value = "  -0.41%  "

# Row 11 - Avalanche
$ws.Range("D11").Value = "'53.71"
$ws.Range("E11").Value = "  +0.97%  "

# Row 12 - ShibaInu
$ws.Range("D12").Value = "'0.0000304"
$ws.Range("E12").Value = "  -0.83%  "

# Row 13 - Polkadot
$ws.Range("D13").Value = "'9.52"
$ws.Range("E13").Value = "  -0.93%  "

# Row 14 - Wrapped liquid staked Ether 2.0
$ws.Range("D14").Value = "'4.086.85"
$ws.Range("E14").Value = "  +0.92%  "

# Row 15 - Bitcoin Cash
$ws.Range("D15").Value = "'597.24"
$ws.Range("E15").Value = "  -1.56%  "

# Row 16 - Wrapped BTC
$ws.Range("D16").Value = "'69.975.16"
$ws.Range("E16").Value = "  +0.81%  "

# Row 17 - Chainlink
$ws.Range("E17").Value = "  +1.28%  "

# Row 18 - Uniswap
$ws.Range("E18").Value = "  +0.51%  "

# Row 19 - Wrapped Ether
$ws.Range("D19").Value = "'3.517.51"
$ws.Range("E19").Value = "  +0.34%  "

# Row 20 - TRON
$ws.Range("E20").Value = "  +1.14%  "

# Row 21 - Polygon
$ws.Range("D21").Value = "'0.993"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22 - Internet Computer (DFINITY)
$ws.Range("D22").Value = "'18.44"
$ws.Range("E22").Value = "  +6.26%  "

# Row 23 - Toncoin
$ws.Range("D23").Value = "'5.27"
$ws.Range("E23").Value = "  +4.92%  "

# Row 24 - Litecoin
$ws.Range("D24").Value = "'102.04"
$ws.Range("E24").Value = "  -2.52%  "

# Row 25 - PancakeSwap
$ws.Range("E25").Value = "  -0.35%  "

# Row 26 - ImmutableX
$ws.Range("D26").Value = "'3.18"
$ws.Range("E26").Value = "  +4.81%  "

# Row 27 - Render Token
$ws.Range("D27").Value = "'10.91"
$ws.Range("E27").Value = "  -0.49%  "

# Row 28 - Filecoin
$ws.Range("E28").Value = "  -3.27%  "

# Row 29 - Ethereum Classic
$ws.Range("D29").Value = "'33.51"
$ws.Range("E29").Value = "  -0.62%  "

# Row 30 - NEAR Protocol
$ws.Range("E30").Value = "  +1.21%  "

# Row 31 - dogwifhat
$ws.Range("D31").Value = "'4.31"
$ws.Range("E31").Value = "  +10.34%  "

# Row 32 - Cosmos
$ws.Range("D32").Value = "'12.51"
$ws.Range("E32").Value = "  +0.09%  "

# Row 33 - Hedera
$ws.Range("E33").Value = "  -1.78%  "

# Row 34 - OKB
$ws.Range("D34").Value = "'63.11"
$ws.Range("E34").Value = "  -0.35%  "

# Row 35 - PEPE
$ws.Range("D35").Value = "'0.0₃0859"
$ws.Range("E35").Value = "  +10.95%  "

# Row 36 - Maker
$ws.Range("D36").Value = "'3.731.67"
$ws.Range("E36").Value = "  +4.55%  "

# Row 37 - Dai
$ws.Range("E37").Value = "  +0.19%  "

# Row 38 - Fetch.AI
$ws.Range("D38").Value = "'3.07"
$ws.Range("E38").Value = "  -3.00%  "

# Row 39 - Stacks
$ws.Range("E39").Value = "  +1.25%  "

# Row 40 - The Graph
$ws.Range("D40").Value = "'0.392"
$ws.Range("E40").Value = "  -0.91%  "

# Row 41 - Injective Protocol
$ws.Range("D41").Value = "'36.61"
$ws.Range("E41").Value = "  -0.25%  "

# Row 42 - Bittensor
$ws.Range("D42").Value = "'489.80"
$ws.Range("E42").Value = "  -6.44%  "

# Row 43 - Kaspa
$ws.Range("E43").Value = "  -3.07%  "

# Row 44 - VeChain
$ws.Range("E44").Value = "  -1.38%  "

# Row 45 - Stellar
$ws.Range("D45").Value = "'0.140"
$ws.Range("E45").Value = "  -3.21%  "

# Row 46 - Theta Token
$ws.Range("E46").Value = "  -4.72%  "

# Row 47 - ApeX Protocol
$ws.Range("D47").Value = "'3.31"
$ws.Range("E47").Value = "  -0.68%  "

# Row 48 - First Digital USD
$ws.Range("D48").Value = "'1.01"
$ws.Range("E48").Value = "  +0.33%  "

# Row 49 - THORChain
$ws.Range("D49").Value = "'8.55"
$ws.Range("E49").Value = "  -3.57%  "

# Row 50 - FLOKI
$ws.Range("D50").Value = "'0.000253"
$ws.Range("E50").Value = "  +4.62%  "

# Row 51 - Monero
$ws.Range("D51").Value = "'130.84"
$ws.Range("E51").Value = "  -0.29%  "
